$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")

# Header row (row 1)
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "autoExtend"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "completion"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "priority"
$ws.Range("G1").Value = "deal"
$ws.Range("H1").Value = "case"
$ws.Range("I1").Value = "tags"
$ws.Range("J1").Value = "description"
$ws.Range("K1").Value = "ownerAssignedTo"
$ws.Range("L1").Value = "keyContact"
$ws.Range("M1").Value = "keyCompany"
$ws.Range("N1").Value = "identifier"

# Task 1 (row 2)
$ws.Range("A2").Value = "Test Title - Task 1"
$ws.Range("B2").Value = "Extend deadline by 14 days"
$ws.Range("C2").Value = "Open"
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = "Meeting"
$ws.Range("F2").Value = "High"
$ws.Range("G2").Value = "Test deal -1"
$ws.Range("H2").Value = "Test case-1"
$ws.Range("I2").Value = "Test tags -1 "
$ws.Range("J2").Value = "Test desc -1 "
$ws.Range("K2").Value = "Tejas niturkar"
$ws.Range("L2").Value = "Test contact - 1"
$ws.Range("M2").Value = "Test company - 1"
$ws.Range("N2").Value = "test identifier - 1"

# Task 2 (row 3)
$ws.Range("A3").Value = "Test Title - Task 2"
$ws.Range("B3").Value = "Extend deadline by 30 days"
$ws.Range("C3").Value = "Open"
$ws.Range("D3").Value = 90
$ws.Range("E3").Value = "Training"
$ws.Range("F3").Value = "Normal"
$ws.Range("G3").Value = "Test deal -2"
$ws.Range("H3").Value = "Test case-2"
$ws.Range("I3").Value = "Test tags -12"
$ws.Range("J3").Value = "Test desc -2 "
$ws.Range("K3").Value = "Tejas niturkar"
$ws.Range("L3").Value = "Test contact - 2"
$ws.Range("M3").Value = "Test company - 2"
$ws.Range("N3").Value = "test identifier - 2"

# completion column is stored as text-formatted numbers (matches source sheet's style)
$ws.Range("D1:D3").NumberFormat = "@"

# Make "tasks" the active sheet/tab and set its selection, like the deals sheet used to be
$ws.Activate()
[void]$ws.Range("M9").Select()
